$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.160.89"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "'2.422.07"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'554.86"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "'137.54"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = "  +3.09%  "
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").Value = "'25.08"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'2.858.34"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "'60.094.15"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "'0.0000138"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "'2.426.80"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "'11.29"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "'4.47"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "'328.14"
$ws.Range("E20").Value = "  -2.58%  "
$ws.Range("D21").Value = "'6.70"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'66.20"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").Value = "'0.178"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").Value = "'8.64"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("D28").Value = "'0.0₃0781"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").Value = "'170.19"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "'6.12"
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("D32").Value = "'1.11"
$ws.Range("E32").Value = "  +6.91%  "
$ws.Range("D33").Value = "'0.405"
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("D34").Value = "'18.56"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.33"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.24"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'330.89"
$ws.Range("E39").Value = "  +3.79%  "
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("D41").Value = "'143.53"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").Value = "'3.68"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").Value = "'0.0970"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "'20.02"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").Value = "'0.0518"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").Value = "'0.576"
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'11.04"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'1.58"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("B50").Value = "ZEEBU"
$ws.Range("C50").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D50").Value = "'4.66"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").Value = "'0.946"
$ws.Range("E51").Value = "  -0.56%  "
